# Extra testfall - problem med avrundning till heltal
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 13: test case 11 (decimal rounding issue with triangle sides)
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Giltig triangel med sidor där decimaltalet börjar på samma heltal"
$ws.Range("C13").Value = "1,2 1,3 1,4"
$ws.Range("D13").Value = "Triangeln har inga lika sidor"
$ws.Range("E13").Value = "Fail"
$ws.Range("F13").Value = "Triangeln är liksidig"

# Add new row 18: extra comment about the rounding bug
$ws.Range("A18").Value = "Applikationen avrundar double värden till heltal, därmed blir resultatet ej korrekt (testfall 11)"

# Widen column B to fit the new, longer text (best-fit recalculated by Excel
# after the longer B13 text was added)
$ws.Columns("B").ColumnWidth = 58.666666666666664

# Move the active selection to the newly added A18 cell
$ws.Range("A18").Select() | Out-Null
